$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 24.02.2022 10:15"

# Row 6 price refresh: new price moves to B6 (current), old price moves to C6 (previous)
$ws.Range("B6").Value = 37.9
$ws.Range("C6").Value = 38.29

# Delta and "old date" columns become plain text values (no longer numeric)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "-0.39"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2022-02-24 10:17:28"
$ws.Range("E6").Style = "Normal"
